$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new column before column C (shifts Timesteps -> D, SimAmount -> E)
$ws.Columns.Item(3).Insert()

# New column inherits column B's width (matches Excel's insert-column behavior)
$ws.Range("C1").ColumnWidth = $ws.Range("B1").ColumnWidth

# Fill in the new column C with header and value
$ws.Range("C1").Value = "Correlation"

# Copy the date cell's format into C2 (reuses the existing date style) then set its value
$ws.Range("A2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "correlationmatrix"

# Update SimAmount value 100 -> 1000 (now in column E)
$ws.Range("E2").Value = 1000

# Update selection to reflect the new active cell
$ws.Range("E3").Select()
